# Update row 6 ("RF") and row 7 ("Ensemble") metric values in the
# evaluation_metrics sheet to reflect the newly recomputed figures
# (added std-dev of precision/recall + PROMISE requirements dummy
# data set changed the downstream classification metrics).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: RF ---
$ws.Range("B6").Value  = 0.8069778082818578
$ws.Range("C6").Value  = 0.8168853728892534
$ws.Range("D6").Value  = 0.8069778082818578
$ws.Range("E6").Value  = 0.8085302072762568
$ws.Range("F6").Value  = 0.8240448409974835
$ws.Range("G6").Value  = 0.8324282209600901
$ws.Range("H6").Value  = 0.8240448409974835
$ws.Range("I6").Value  = 0.8256457789346012
$ws.Range("J6").Value  = 0.7532601235415237
$ws.Range("K6").Value  = 0.7594462895492863
$ws.Range("L6").Value  = 0.7532601235415237
$ws.Range("M6").Value  = 0.7513773768325669
$ws.Range("N6").Value  = 0.8090597117364448
$ws.Range("O6").Value  = 0.8202125323898073
$ws.Range("P6").Value  = 0.8090597117364448
$ws.Range("Q6").Value  = 0.8095402314725835
$ws.Range("R6").Value  = 0.8005033173186915
$ws.Range("S6").Value  = 0.8162938621815732
$ws.Range("T6").Value  = 0.8005033173186915
$ws.Range("U6").Value  = 0.8030358157294298
$ws.Range("V6").Value  = 0.8198123999084878
$ws.Range("W6").Value  = 0.8243610942814212
$ws.Range("X6").Value  = 0.8198123999084878
$ws.Range("Y6").Value  = 0.8189970648453558

# --- Row 7: Ensemble ---
$ws.Range("B7").Value  = 0.8455044612216884
$ws.Range("C7").Value  = 0.8503641615949924
$ws.Range("D7").Value  = 0.8455044612216884
$ws.Range("E7").Value  = 0.8460971644411467
$ws.Range("F7").Value  = 0.8627316403568978
$ws.Range("G7").Value  = 0.8650388708177061
$ws.Range("H7").Value  = 0.8627316403568978
$ws.Range("I7").Value  = 0.8618252690317654
$ws.Range("J7").Value  = 0.8305422100205903
$ws.Range("K7").Value  = 0.8424118199613273
$ws.Range("L7").Value  = 0.8305422100205903
$ws.Range("M7").Value  = 0.8313807466572986
$ws.Range("N7").Value  = 0.8498055364905056
$ws.Range("O7").Value  = 0.8549274283698096
$ws.Range("P7").Value  = 0.8498055364905056
$ws.Range("Q7").Value  = 0.8499449907477221
$ws.Range("R7").Value  = 0.8690459849004805
$ws.Range("S7").Value  = 0.8754127222478628
$ws.Range("T7").Value  = 0.8690459849004805
$ws.Range("U7").Value  = 0.8689077229551714
$ws.Range("V7").Value  = 0.8347517730496454
$ws.Range("W7").Value  = 0.8404846004677561
$ws.Range("X7").Value  = 0.8347517730496454
$ws.Range("Y7").Value  = 0.8351239298215546
